# Semaine 9 update (v9.0): fill in the previously-empty bullet paragraph
# at the end of the "Semaine 9" section with two new list items, reusing
# the same numbering definition (numId 1) already used by the other
# "Paragraphedeliste" bullets in the document.

$d = $word.ActiveDocument

# The last paragraph in the document is the empty "Paragraphedeliste"
# placeholder sitting right after the "Semaine 9" / "18/01 - 22/01" lines.
$count = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($count)

# Grab an existing numbered-list template (any "Paragraphedeliste" bullet
# elsewhere in the doc carries it) so we continue list numId 1 instead of
# minting a brand new numbering instance.
$listTemplate = $d.Paragraphs.Item(102).Range.ListFormat.ListTemplate

# First new bullet: turn the empty paragraph into a real list item.
$target.Range.Text = "Retour sur OFv4"
$target.Range.ListFormat.ApplyListTemplate($listTemplate, $true)

# Second new bullet: insert a fresh paragraph after it, still on the list.
$newPara = $target.Range.InsertParagraphAfter()
$d.Paragraphs.Item($d.Paragraphs.Count).Range.Text = "Il y a eu une migration de SVN vers Git, j'ai donc importer le `"nouveau`" projet avec git, "
